$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 104.629996108622
$ws.Range("B2").Value = 0.1095773622483572
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 123.772948660898
$ws.Range("B3").Value = -0.05333307395443487
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 96.47495660750897
$ws.Range("B4").Value = 0.1789784140845004
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 106.471634287554
$ws.Range("B5").Value = 0.1187283637482496
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 104.752381984892
$ws.Range("B6").Value = 0.13295871063862
$ws.Range("C6").Value = 5
